# Apply updates to column F ("dSF") values per the repull/mean recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -3
$ws.Range("F7").Value = 0
$ws.Range("F11").Value = -7
$ws.Range("F16").Value = -4
$ws.Range("F17").Value = 8
$ws.Range("F18").Value = -4
$ws.Range("F23").Value = -5
